$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "26.275.98"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.624.13"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "212.68"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  +0.04%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.249"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.76%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.0615"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "18.94"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +4.63%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0816"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.61%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "1.850.86"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.627.54"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +1.64%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "4.03"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +1.30%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "26.290.80"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  +1.11%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "62.50"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +3.85%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0729"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("E19").Value = "  +0.07%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "203.31"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.81%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.29"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("E22").Value = "  +1.10%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "6.04"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.55%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +7.82%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "142.76"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +0.03%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "15.26"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "6.56"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0526"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +10.57%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  +2.20%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.40"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.79%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.171.07"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("E37").Value = "  +1.49%  "
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("E39").Value = "  +0.07%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "2.32"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.23%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.498"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("E42").Value = "  +1.08%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.31"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.39%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.762.16"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +14.34%  "
$ws.Range("E47").Value = "  +0.86%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "54.20"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.18%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0508"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  +0.01%  "

Write-Host "Updated cryptos list values."
